# Updated cryptos list on Tue Nov 19 15:45:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text that can look
# like numbers/percentages (e.g. "614.15", "0.0000251"). Force the whole
# data range to Text format first so assigning those strings doesn't get
# auto-coerced into numeric cells (which would also mangle formatting via
# floating point rounding). The format + style are reset back afterwards
# so the cells end up with the same (default) style they started with.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "92.172.76"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.103.60"
$ws.Range("E3").Value = "  -1.10%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.21%  "

# Row 5 - Solana
$ws.Range("D5").Value = "239.78"
$ws.Range("E5").Value = "  -1.62%  "

# Row 6 - BNB
$ws.Range("D6").Value = "614.15"
$ws.Range("E6").Value = "  -1.46%  "

# Row 7 - XRP
$ws.Range("D7").Value = "1.10"
$ws.Range("E7").Value = "  -5.75%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "0.392"
$ws.Range("E8").Value = "  +5.16%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - LidoStakedEther
$ws.Range("D10").Value = "3.100.65"
$ws.Range("E10").Value = "  -1.05%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.89%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.64%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -0.14%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "34.33"
$ws.Range("E14").Value = "  -3.60%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "91.748.31"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16 - Toncoin
$ws.Range("E16").Value = "  -0.40%  "

# Row 17 - WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "3.675.35"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.087.96"
$ws.Range("E18").Value = "  -1.81%  "

# Row 19 - SuiNetwork
$ws.Range("D19").Value = "3.66"
$ws.Range("E19").Value = "  -2.78%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").Value = "  -1.97%  "

# Row 21 - Polkadot
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  -0.46%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "447.62"
$ws.Range("E22").Value = "  +0.47%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +0.84%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  -5.11%  "

# Row 25 - NEARProtocol
$ws.Range("E25").Value = "  -1.58%  "

# Row 26 - now Litecoin (was Aptos)
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "87.10"
$ws.Range("E26").Value = "  -3.60%  "

# Row 27 - now Aptos (was Litecoin)
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "11.72"
$ws.Range("E27").Value = "  -2.14%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "3.265.84"
$ws.Range("E28").Value = "  +0.60%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.14%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +12.30%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.237"
$ws.Range("E31").Value = "  -5.40%  "

# Row 32 - Cronos
$ws.Range("E32").Value = "  -4.85%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "9.17"
$ws.Range("E33").Value = "  -1.54%  "

# Row 34 - Binance-PegBSC-USD
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +57.17%  "

# Row 35 - RenderToken (only price changed, volume unchanged)
$ws.Range("D35").Value = "7.96"

# Row 36 - Kaspa
$ws.Range("D36").Value = "0.162"
$ws.Range("E36").Value = "  -4.95%  "

# Row 37 - MantraDAO
$ws.Range("D37").Value = "4.24"
$ws.Range("E37").Value = "  -3.10%  "

# Row 38 - EthereumClassic
$ws.Range("D38").Value = "26.16"
$ws.Range("E38").Value = "  -2.12%  "

# Row 39 - PancakeSwap
$ws.Range("E39").Value = "  +0.01%  "

# Row 40 - Bittensor
$ws.Range("E40").Value = "  -3.52%  "

# Row 42 - PolygonEcosystemToken
$ws.Range("E42").Value = "  +1.80%  "

# Row 43 - dogwifhat
$ws.Range("E43").Value = "  -4.41%  "

# Row 44 - WhiteBITCoin
$ws.Range("D44").Value = "22.44"
$ws.Range("E44").Value = "  +0.86%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.02%  "

# Row 46 - Monero
$ws.Range("D46").Value = "159.54"
$ws.Range("E46").Value = "  +3.54%  "

# Row 47 - Stacks
$ws.Range("D47").Value = "1.89"
$ws.Range("E47").Value = "  -2.59%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  -1.58%  "

# Row 49 - ImmutableX
$ws.Range("E49").Value = "  +0.85%  "

# Row 50 - VeChain
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +3.91%  "

# Row 51 - OKB
$ws.Range("D51").Value = "44.10"
$ws.Range("E51").Value = "  -0.91%  "

# Restore the original (default) style now that the text values are set,
# so no cell ends up with a lingering explicit number format.
$dataRange.Style = "Normal"
